$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "settings" sheet: the form_id column is being retired (pyxform/cht-conf
#    no longer uses it). Before removing the column, roll the header-cell
#    comments down one slot so the text that used to document "version"
#    now documents "form_id"'s old slot (B1), "pages" documents what used to
#    be "version"'s slot (C1), and "namespaces" documents the comment that
#    used to sit on "pages" (D1) -- then drop the now-duplicate trailing
#    comment that used to live on "namespaces" (E1).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("settings")
$ws2.Activate()

$versionText    = $ws2.Range("C1").Comment.Text()
$pagesText      = $ws2.Range("D1").Comment.Text()
$namespacesText = $ws2.Range("E1").Comment.Text()

$ws2.Range("B1").Comment.Text($versionText)
$ws2.Range("C1").Comment.Text($pagesText)
$ws2.Range("D1").Comment.Text($namespacesText)
$ws2.Range("E1").Comment.Delete()

# Now actually remove the form_id column (B): title row + value row + all
# comments/column widths shift left automatically.
$ws2.Columns.Item(2).Delete()

$ws2.Range("B1").Select()

# ---------------------------------------------------------------------------
# 2. "survey" sheet: tidy up the conditional formatting so the three blocks
#    each cover a single uniform range instead of the old patchwork that
#    carved out a hole around row 27, and move the cursor to A6.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("survey")
$ws1.Activate()

# Drop the two rule-sets that existed only to cover the C27 gap.
$ws1.Range("C27").FormatConditions.Delete()

# Consolidate the 5-rule block (currently "A28:D10000 A27:B27 D2:D27 A2:C26")
# down to a single "A2:D10000" block, keeping the same dxf-backed rules.
$fcs = $ws1.Range("C2").FormatConditions
$fcs.Item(1).ModifyAppliesToRange($ws1.Range("A2:D10000"))

# Consolidate the C-column block (currently "C28:C10000 C2:C26") down to
# "C2:C10000".
$fcs.Item($fcs.Count).ModifyAppliesToRange($ws1.Range("C2:C10000"))

$ws1.Range("A6").Select()
